# This script re-applies an "automatic update of files" style refresh to the
# Artfynd export sheet: the underlying data rows (2-9, 11, 12, 14) get new
# observation records (new Id, possibly new species/taxon, and new
# coordinates), while rows 1 (header), 10 and 13 are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111523701
$ws.Range("B2").Value = 89686
$ws.Range("D2").Value = 'NT'
$ws.Range("E2").Value = 658
$ws.Range("F2").Value = 'Rosenticka'
$ws.Range("G2").Value = 'Rhodofomes roseus'
$ws.Range("H2").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q2").Value = 497367.2942720717
$ws.Range("R2").Value = 6754083.757028132

# Row 3
$ws.Range("A3").Value = 111523695
$ws.Range("B3").Value = 5113
$ws.Range("D3").Value = 'LC'
$ws.Range("E3").Value = 100526
$ws.Range("F3").Value = 'Bronshjon'
$ws.Range("G3").Value = 'Callidium coriaceum'
$ws.Range("H3").Value = 'Paykull, 1800'
$ws.Range("Q3").Value = 497354.1644349985
$ws.Range("R3").Value = 6754111.484663551

# Row 4
$ws.Range("A4").Value = 111523683
$ws.Range("B4").Value = 89845
$ws.Range("D4").Value = 'VU'
$ws.Range("E4").Value = 1209
$ws.Range("F4").Value = 'Rynkskinn'
$ws.Range("G4").Value = 'Phlebia centrifuga'
$ws.Range("H4").Value = 'P.Karst.'
$ws.Range("Q4").Value = 497391.6869587752
$ws.Range("R4").Value = 6754138.20205555

# Row 5
$ws.Range("A5").Value = 111523724
$ws.Range("B5").Value = 93881
$ws.Range("D5").Value = 'LC'
$ws.Range("E5").Value = 2869
$ws.Range("F5").Value = 'Bollvitmossa'
$ws.Range("G5").Value = 'Sphagnum wulfianum'
$ws.Range("H5").Value = 'Girg.'
$ws.Range("Q5").Value = 497291.3182300103
$ws.Range("R5").Value = 6754089.649475355

# Row 6
$ws.Range("A6").Value = 111523657
$ws.Range("B6").Value = 89686
$ws.Range("D6").Value = 'NT'
$ws.Range("E6").Value = 658
$ws.Range("F6").Value = 'Rosenticka'
$ws.Range("G6").Value = 'Rhodofomes roseus'
$ws.Range("H6").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q6").Value = 497390.1961838813
$ws.Range("R6").Value = 6754097.842248607

# Row 7
$ws.Range("A7").Value = 111523731
$ws.Range("B7").Value = 89845
$ws.Range("D7").Value = 'VU'
$ws.Range("E7").Value = 1209
$ws.Range("F7").Value = 'Rynkskinn'
$ws.Range("G7").Value = 'Phlebia centrifuga'
$ws.Range("H7").Value = 'P.Karst.'
$ws.Range("Q7").Value = 497307.3714758331
$ws.Range("R7").Value = 6754063.864355386

# Row 8
$ws.Range("A8").Value = 111523727
$ws.Range("B8").Value = 89845
$ws.Range("D8").Value = 'VU'
$ws.Range("E8").Value = 1209
$ws.Range("F8").Value = 'Rynkskinn'
$ws.Range("G8").Value = 'Phlebia centrifuga'
$ws.Range("H8").Value = 'P.Karst.'
$ws.Range("Q8").Value = 497338.5868253836
$ws.Range("R8").Value = 6754122.194367126

# Row 9
$ws.Range("A9").Value = 111523656
$ws.Range("B9").Value = 89845
$ws.Range("D9").Value = 'VU'
$ws.Range("E9").Value = 1209
$ws.Range("F9").Value = 'Rynkskinn'
$ws.Range("G9").Value = 'Phlebia centrifuga'
$ws.Range("H9").Value = 'P.Karst.'
$ws.Range("Q9").Value = 497390.1961838813
$ws.Range("R9").Value = 6754097.842248607

# Row 11
$ws.Range("A11").Value = 111523741
$ws.Range("B11").Value = 89686
$ws.Range("D11").Value = 'NT'
$ws.Range("E11").Value = 658
$ws.Range("F11").Value = 'Rosenticka'
$ws.Range("G11").Value = 'Rhodofomes roseus'
$ws.Range("H11").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q11").Value = 497384.3941364431
$ws.Range("R11").Value = 6754155.713205664

# Row 12
$ws.Range("A12").Value = 111523728
$ws.Range("B12").Value = 89686
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 658
$ws.Range("F12").Value = 'Rosenticka'
$ws.Range("G12").Value = 'Rhodofomes roseus'
$ws.Range("H12").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q12").Value = 497338.5868253836
$ws.Range("R12").Value = 6754122.194367126

# Row 14
$ws.Range("A14").Value = 111523730
$ws.Range("B14").Value = 89405
$ws.Range("D14").Value = 'NT'
$ws.Range("E14").Value = 1202
$ws.Range("F14").Value = 'Ullticka'
$ws.Range("G14").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H14").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q14").Value = 497338.5868253836
$ws.Range("R14").Value = 6754122.194367126
